# "Mise a jour apogee d'avril 2018"
#
# On the "Situation_sociale" sheet, insert a new code/label pair
#   B6 = "DD"   C6 = "Demi Droit"
# right after the existing "CH / Chomeur" row (row 5), pushing the
# rows that used to be 6-10 ("NO/Normal" ... "TH/Thesard") down to 7-11,
# and the trailing blank spacer row from 11 to 12.
#
# This also implicitly grows the shared-strings table with the two new
# strings "DD" and "Demi Droit".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Situation_sociale")

# xlShiftDown = -4121 ; insert a brand-new blank row at position 6, which
# shifts the current rows 6-10 down to 7-11 (values & formats travel
# together with their rows).
$ws.Rows.Item(6).Insert(-4121)

# xlPasteFormats = -4122 ; clone the visual formatting (fill / border /
# font / number format / alignment) of a normal data row onto the new
# row 6, so it matches the style used by the "even" rows of the table
# (same style as row 4 / BO-Boursier).
$ws.Range("B4:C4").Copy()
$ws.Range("B6:C6").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Rows.Item(6).RowHeight = 19.7321

# Fill in the new row's values.
$ws.Cells.Item(6, 2).Value = "DD"
$ws.Cells.Item(6, 3).Value = "Demi Droit"

# The banded (alternating) row styling is tied to the absolute row
# number, not to the data that travelled down with Insert(), so re-sync
# rows 7-11 to the correct alternating pattern: 7/9/11 use the "odd"
# style (same as row 5 / CH-Chomeur) and 8/10 use the "even" style
# (same as row 4 / BO-Boursier).
$ws.Range("B5:C5").Copy()
$ws.Range("B7:C7").PasteSpecial(-4122)
$ws.Range("B9:C9").PasteSpecial(-4122)
$ws.Range("B11:C11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B4:C4").Copy()
$ws.Range("B8:C8").PasteSpecial(-4122)
$ws.Range("B10:C10").PasteSpecial(-4122)
$excel.CutCopyMode = 0
